$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 633
$ws.Range("F5").Value = 2795
$ws.Range("F9").Value = 6465
$ws.Range("F13").Value = 5051
$ws.Range("F16").Value = 557
$ws.Range("F19").Value = 1525
$ws.Range("F22").Value = 126
$ws.Range("F23").Value = 135
$ws.Range("F24").Value = 1088
$ws.Range("F25").Value = 251
$ws.Range("F26").Value = 1395
$ws.Range("F27").Value = 1059
$ws.Range("F29").Value = 331
$ws.Range("F30").Value = 593
$ws.Range("F34").Value = 257
$ws.Range("F35").Value = 1523
$ws.Range("F38").Value = 618
$ws.Range("F39").Value = 1081
$ws.Range("F40").Value = 121
$ws.Range("F41").Value = 553
$ws.Range("F43").Value = 2315
$ws.Range("F44").Value = 2576
$ws.Range("F46").Value = 143
$ws.Range("F49").Value = 390

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G5").Value = 230
$ws.Range("F9").Value = 328
$ws.Range("F12").Value = 96
$ws.Range("F16").Value = 241
$ws.Range("F17").Value = 162
$ws.Range("F21").Value = 154
$ws.Range("F30").Value = 17
$ws.Range("F41").Value = 31

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 517
$ws.Range("F8").Value = 1533
$ws.Range("F10").Value = 2575
$ws.Range("F11").Value = 881
$ws.Range("F12").Value = 756
$ws.Range("F14").Value = 138

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 517
$ws.Range("F5").Value = 633
$ws.Range("F7").Value = 2795
$ws.Range("F9").Value = 1533
$ws.Range("F11").Value = 2575
$ws.Range("F12").Value = 6465
$ws.Range("F13").Value = 756
$ws.Range("F17").Value = 1525
$ws.Range("F19").Value = 126
$ws.Range("F20").Value = 135
$ws.Range("F21").Value = 1088
$ws.Range("F22").Value = 251
$ws.Range("F23").Value = 96
$ws.Range("F24").Value = 1395
$ws.Range("F25").Value = 1059
$ws.Range("F27").Value = 331
$ws.Range("F28").Value = 593
$ws.Range("F31").Value = 257
$ws.Range("F33").Value = 1523
$ws.Range("F36").Value = 1081
$ws.Range("F37").Value = 553
$ws.Range("F43").Value = 2315
$ws.Range("F44").Value = 2576
$ws.Range("F48").Value = 390

